# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect refreshed scrape data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 83
    $ws.Range("F4").Value = 2226
    $ws.Range("F5").Value = 197
}
